$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.638.45'
$ws.Range("E2").Value = '  +3.74%  '
$ws.Range("D3").Value = '1.915.11'
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.09'
$ws.Range("E5").Value = '  +1.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.701'
$ws.Range("E6").Value = '  +2.98%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.87'
$ws.Range("E8").Value = '  +2.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.370'
$ws.Range("E9").Value = '  +3.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.44'
$ws.Range("E10").Value = '  +9.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0764'
$ws.Range("E11").Value = '  +3.01%  '
$ws.Range("E12").Value = '  +2.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.58'
$ws.Range("E13").Value = '  +8.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.814'
$ws.Range("E14").Value = '  +6.22%  '
$ws.Range("D15").Value = '2.191.79'
$ws.Range("E15").Value = '  +1.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.13'
$ws.Range("E16").Value = '  +3.99%  '
$ws.Range("D17").Value = '1.914.95'
$ws.Range("E17").Value = '  +1.93%  '
$ws.Range("D18").Value = '36.583.00'
$ws.Range("E18").Value = '  +3.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.56'
$ws.Range("E19").Value = '  +1.41%  '
$ws.Range("D20").Value = '0.0₃0864'
$ws.Range("E20").Value = '  +4.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '250.49'
$ws.Range("E21").Value = '  +2.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.41'
$ws.Range("E22").Value = '  +4.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.22'
$ws.Range("E23").Value = '  +3.58%  '
$ws.Range("E24").Value = '  -2.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  +2.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '168.90'
$ws.Range("E27").Value = '  +2.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.79'
$ws.Range("E28").Value = '  +1.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.74'
$ws.Range("E29").Value = '  +2.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.130'
$ws.Range("E30").Value = '  +1.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.57'
$ws.Range("E31").Value = '  +6.31%  '
$ws.Range("E32").Value = '  +4.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.35'
$ws.Range("E33").Value = '  +4.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0892'
$ws.Range("E34").Value = '  +22.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.90'
$ws.Range("E35").Value = '  +3.42%  '
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("E37").Value = '  +6.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.886'
$ws.Range("E38").Value = '  +3.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.79'
$ws.Range("E39").Value = '  +50.32%  '
$ws.Range("E40").Value = '  +4.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '106.55'
$ws.Range("E41").Value = '  +10.03%  '
$ws.Range("E42").Value = '  +3.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.60'
$ws.Range("E43").Value = '  +1.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.96'
$ws.Range("E44").Value = '  +23.83%  '
$ws.Range("E45").Value = '  +3.79%  '
$ws.Range("D46").Value = '1.345.24'
$ws.Range("E46").Value = '  +2.95%  '
$ws.Range("E47").Value = '  -1.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0813'
$ws.Range("E48").Value = '  +1.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.80'
$ws.Range("E49").Value = '  +2.66%  '
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.71'
$ws.Range("E50").Value = '  +4.15%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.40'
$ws.Range("E51").Value = '  +1.61%  '
